$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Diclofenac")

# Update the DrugBank row's "Code" cell: was text "C0012091", now a plain
# number 3355.
$ws.Range("D2").Value = 3355

# Move/record the current selection (matches the saved sheet view state).
$null = $ws.Range("E11").Select()
